$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 0.27610753363967433
$ws.Range("A2").Value = -0.009999999279230565
$ws.Range("A3").Value = -0.0089999992768508
$ws.Range("A4").Value = 0.061999999502965863
$ws.Range("A5").Value = -0.005999999294264313
$ws.Range("A6").Value = -0.005999999279225676
$ws.Range("A7").Value = -0.01999999915619277
$ws.Range("A8").Value = -0.019999999148916814
$ws.Range("A9").Value = -0.030400446508166645
$ws.Range("A10").Value = -0.005999999257738864
$ws.Range("A11").Value = -0.004499999270574762
$ws.Range("A12").Value = -0.005999999257414235
$ws.Range("A13").Value = -0.005999999256588673
$ws.Range("A14").Value = -0.011999999205175804
$ws.Range("A15").Value = 0.03475219976130539
$ws.Range("A16").Value = -0.005999999255796862
$ws.Range("A17").Value = -0.00599999925309902
$ws.Range("A18").Value = -0.00899999922677619
$ws.Range("A19").Value = -0.08378886869011293
$ws.Range("A20").Value = -0.008999999279485138
$ws.Range("A21").Value = -0.008999999278681337
$ws.Range("A22").Value = -0.008999999278130666
$ws.Range("A23").Value = -0.008999999270657533
$ws.Range("A24").Value = -0.041999998975218666
$ws.Range("A25").Value = -0.04199999896990736
$ws.Range("A26").Value = -0.0059999992791262
$ws.Range("A27").Value = -0.005999999279008961
$ws.Range("A28").Value = 0.06714561575174471
$ws.Range("A29").Value = -0.011999999219165502
$ws.Range("A30").Value = -0.00796672798882847
$ws.Range("A31").Value = -0.014999999182936818
$ws.Range("A32").Value = -0.02099999913024142
$ws.Range("A33").Value = -0.005999999259167055
